{"js": "// 1. Update activation date: 01/01/2012 -> 01/01/2023\nconst body = context.document.body;\nconst dateResults = body.search(\"01/01/2012\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < dateResults.items.length; i++) {\n  dateResults.items[i].insertText(\"01/01/2023\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2-4. Insert italic English-translation paragraphs right after the\n// matching Portuguese paragraphs.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst insertions = [\n  {\n    marker: \"A reologia \u00e9 a ci\u00eancia\",\n    text:\n      \"Rheology is the science that studies the flow of materials. Your knowledge is necessary to understand the processes of forming materials. The course aims to provide students with the basic and applied concepts of rheology and familiarize them with experimental methods for evaluating the rheological properties of materials.\",\n  },\n  {\n    marker: \"Escoamento de fluidos newtonianos\",\n    text:\n      \"Flow of Newtonian and non-Newtonian fluids. Viscosity and rheometry. viscoelasticity. Applications.\",\n  },\n  {\n    marker: \"1. Introdu\u00e7\u00e3o.\",\n    text:\n      \"1. Introduction. 2. Stress and deformation. 3. Types of deformation and flow of materials. 4. Fundamental equations of rheology. Flow of Newtonian and non-Newtonian fluids. 5. Viscosimetry and rheometry. 6. Rheology of dispersed systems. Colloids and emulsions. diluted solutions. Capillary viscosimetry. 7. Rheology of molten polymers. 8. Viscoelasticity. 9. Dynamic-mechanical behavior of materials. 10. Applications.\",\n  },\n];\n\nfor (const { marker, text } of insertions) {\n  const target = paragraphs.items.find((p) => p.text.startsWith(marker));\n  const newPara = target.insertParagraph(text, Word.InsertLocation.after);\n  newPara.font.italic = true;\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update activation date\n$find = $d.Content.Find\n$find.Text = \"Ativa\u00e7\u00e3o: 01/01/2012\"\n$find.Replacement.Text = \"Ativa\u00e7\u00e3o: 01/01/2023\"\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# Helper: insert a new italic paragraph right after the paragraph whose\n# text starts with $marker, with text $newText.\nfunction Insert-ItalicParagraphAfter($marker, $newText) {\n    $doc = $word.ActiveDocument\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $p = $doc.Paragraphs($i)\n        if ($p.Range.Text.StartsWith($marker)) {\n            $p.Range.InsertParagraphAfter()\n            $newPara = $doc.Paragraphs($i + 1)\n            $r = $newPara.Range\n            $r.Text = $newText\n            $r2 = $doc.Range($r.Start, $r.Start + $newText.Length)\n            $r2.Font.Italic = 1\n            break\n        }\n    }\n}\n\n# Go from the bottom-most insertion point to the top-most so earlier\n# paragraph indices stay valid while we work.\n\n# 4. \"Programa\" body -> add English translation paragraph after it\nInsert-ItalicParagraphAfter \"1. Introdu\u00e7\u00e3o.\" \"1. Introduction. 2. Stress and deformation. 3. Types of deformation and flow of materials. 4. Fundamental equations of rheology. Flow of Newtonian and non-Newtonian fluids. 5. Viscosimetry and rheometry. 6. Rheology of dispersed systems. Colloids and emulsions. diluted solutions. Capillary viscosimetry. 7. Rheology of molten polymers. 8. Viscoelasticity. 9. Dynamic-mechanical behavior of materials. 10. Applications.\"\n\n# 3. \"Programa resumido\" body -> add English translation paragraph after it\nInsert-ItalicParagraphAfter \"Escoamento de fluidos newtonianos\" \"Flow of Newtonian and non-Newtonian fluids. Viscosity and rheometry. viscoelasticity. Applications.\"\n\n# 2. \"Objetivos\" body -> add English translation paragraph after it\nInsert-ItalicParagraphAfter \"A reologia \u00e9 a ci\u00eancia\" \"Rheology is the science that studies the flow of materials. Your knowledge is necessary to understand the processes of forming materials. The course aims to provide students with the basic and applied concepts of rheology and familiarize them with experimental methods for evaluating the rheological properties of materials.\"\n"}
